# Fix E2E tests and ingestion logic - update fixture data to match generated fixtures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "TOOL"
$ws.Range("B1").Value = "TYPE"
$ws.Range("C1").Value = "AREA"
$ws.Range("D1").Value = "STATION"
$ws.Range("E1").Value = "DESCRIPTION"

# Update row 2
$ws.Range("A2").Value = "WeldGun01"
$ws.Range("B2").Value = "Weld Gun"
$ws.Range("C2").Value = "Body Shop"
$ws.Range("D2").Value = "OP10"
$ws.Range("E2").Value = "Servo Gun"

# Update row 3
$ws.Range("A3").Value = "Gripper01"
$ws.Range("B3").Value = "Gripper"
$ws.Range("C3").Value = "Body Shop"
$ws.Range("D3").Value = "OP20"
$ws.Range("E3").Value = "Material Handling"

# Remove row 4 entirely (it existed before, now the data only spans to row 3)
$ws.Range("A4:E4").Delete()
